# Insert a new weekly price observation row above the current row 462
# (i.e. at row 461), shifting every subsequent row down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(461).EntireRow.Insert()

$ws.Range("A461").Value = 4
$ws.Range("B461").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C461").Value = "Los Lagos"
$ws.Range("D461").Value = 45135
$ws.Range("E461").Value = 10
$ws.Range("F461").Value = 100112008
$ws.Range("G461").Value = "Coliflor"
$ws.Range("H461").Value = "Sin especificar"
$ws.Range("I461").Value = "Primera"
$ws.Range("J461").Value = 800
$ws.Range("K461").Value = 1400
$ws.Range("L461").Value = 1400
$ws.Range("M461").Value = 1400
$ws.Range("N461").Value = "$/unidad"
$ws.Range("O461").Value = "Región Metropolitana"
$ws.Range("P461").Value = 1400
$ws.Range("Q461").Value = 1
$ws.Range("R461").Value = "Hortaliza"
